$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UF_IVP_DIARIO")

# Copy the formatting (styles) of the last existing data row down to the new rows
$ws.Range("A833:C833").Copy()
$ws.Range("A834:C863").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @(44296, 29415.74, 30576.66),
    @(44297, 29419.66, 30580.26),
    @(44298, 29423.57, 30583.85),
    @(44299, 29427.49, 30587.45),
    @(44300, 29431.41, 30591.05),
    @(44301, 29435.32, 30594.65),
    @(44302, 29439.24, 30598.25),
    @(44303, 29443.16, 30601.85),
    @(44304, 29447.07, 30605.45),
    @(44305, 29450.99, 30609.05),
    @(44306, 29454.91, 30612.65),
    @(44307, 29458.83, 30616.26),
    @(44308, 29462.75, 30619.86),
    @(44309, 29466.67, 30623.46),
    @(44310, 29470.59, 30627.06),
    @(44311, 29474.52, 30630.67),
    @(44312, 29478.44, 30634.27),
    @(44313, 29482.36, 30637.88),
    @(44314, 29486.29, 30641.48),
    @(44315, 29490.21, 30645.09),
    @(44316, 29494.13, 30648.69),
    @(44317, 29498.06, 30652.3),
    @(44318, 29501.98, 30655.91),
    @(44319, 29505.91, 30659.51),
    @(44320, 29509.84, 30663.119999999999),
    @(44321, 29513.759999999998, 30666.73),
    @(44322, 29517.69, 30670.34),
    @(44323, 29521.62, 30673.95),
    @(44324, 29525.55, 30677.56),
    @(44325, 29529.48, 30681.17),
)

$startRow = 834
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Extend the sheet-scoped defined name that tracks the data range
$definedName = $wb.Names.Item("UF_IVP_DIARIO!UF_IVP_DIARIO")
$definedName.RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$863"

# Update the view to match Excel's behaviour of scrolling to the newly entered data
$ws.Activate()
$ws.Range("A863").Select() | Out-Null
